$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.970.70"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "'1.819.04"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'310.24"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").Value = "'0.4687"
$ws.Range("E7").Value = "  +0.88%  "
$ws.Range("E8").Value = "  -1.15%  "
$ws.Range("D9").Value = "'0.07353"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").Value = "'0.8723"
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("D12").Value = "'1.821.10"
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("D13").Value = "'5.411"
$ws.Range("E13").Value = "  +1.15%  "
$ws.Range("D14").Value = "'0.07106"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").Value = "'6.509"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").Value = "'91.48"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").Value = "'0.000008719"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").Value = "'14.67"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").Value = "'26.985.62"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").Value = "'5.286"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("D24").Value = "'2.043.14"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").Value = "'1.893"
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("D26").Value = "'151.11"
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("D27").Value = "'18.35"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").Value = "'2.153"
$ws.Range("E28").Value = "  +0.60%  "
$ws.Range("D29").Value = "'5.249"
$ws.Range("E29").Value = "  -1.22%  "
$ws.Range("D30").Value = "'116.90"
$ws.Range("E30").Value = "  +1.02%  "
$ws.Range("D31").Value = "'0.08884"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("D32").Value = "'0.7585"
$ws.Range("E32").Value = "  +0.29%  "
$ws.Range("D33").Value = "'1.162"
$ws.Range("E33").Value = "  +0.50%  "
$ws.Range("D34").Value = "'4.509"
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("D35").Value = "'2.928"
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "'0.05307"
$ws.Range("E38").Value = "  +0.98%  "
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("E40").Value = "  +1.50%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'7.182"
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'2.367"
$ws.Range("E42").Value = "  -2.72%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.5288"
$ws.Range("E43").Value = "  -0.82%  "
$ws.Range("D44").Value = "'0.1653"
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("D45").Value = "'8.442"
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("D46").Value = "'0.4873"
$ws.Range("E46").Value = "  -1.73%  "
$ws.Range("E47").Value = "  +1.28%  "
$ws.Range("D48").Value = "'1.002"
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("D49").Value = "'103.41"
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("E50").Value = "  -1.01%  "
$ws.Range("D51").Value = "'0.06298"
$ws.Range("E51").Value = "  +0.12%  "
